# Update the "Weight" (D) and "Percent Change" (E) columns of the QE
# holdings table (rows 2-35) with refreshed values, as uploaded by the
# author. The worksheet is protected (legacy password "D382"), so it has
# to be unprotected before the edit and re-protected with the same
# password afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$wasProtected = $ws.ProtectContents
if ($wasProtected) {
    $ws.Unprotect("D382")
}

# row, new Weight (column D), new Percent Change (column E)
# $null means "leave column D unchanged" (row 35 / Total keeps Weight = 1)
$updates = @(
    @(2, 0.09335656915765851, 0.0179552175749893),
    @(3, 0.07935619904212959, 0.008175466397752684),
    @(4, 0.05252696734623594, 0.005669606664236326),
    @(5, 0.0515183668044954, 0.00692340977931627),
    @(6, 0.04793496255616932, 0.001669449081802998),
    @(7, 0.04165285878520653, 0.00117332638028822),
    @(8, 0.03606554331602079, 0.002250574530469951),
    @(9, 0.03840513074438236, 0.008572959457045837),
    @(10, 0.03399604019254582, -0.001643561526368353),
    @(11, 0.03517612283262704, 0.003412470140886414),
    @(12, 0.03528901521403467, -0.006491297692612807),
    @(13, 0.03092079456230267, -0.007986024457199981),
    @(14, 0.03193596858447958, 0.01321806677942772),
    @(15, 0.03224892341395136, -0.00584919706476672),
    @(16, 0.02978472566079555, -0.00500892394495922),
    @(17, 0.02914766966545986, 0.007501103103397577),
    @(18, 0.02827654060588913, -0.02263561660450608),
    @(19, 0.02389517299328433, 0.008844953173777537),
    @(20, 0.02095882786269714, 0.002127282396738206),
    @(21, 0.02192470078160112, -0.01564282222584323),
    @(22, 0.02155872940339235, 0.008166295471417895),
    @(23, 0.0208399336078476, -0.002317703690497441),
    @(24, 0.01907624022641353, -0.003101309441764344),
    @(25, 0.02137009909521758, 0.03623013962445842),
    @(26, 0.02021259493141788, 0.00152710613387641),
    @(27, 0.01968443006852854, 0.005285013212533052),
    @(28, 0.01867754434777141, 0.008836896145430151),
    @(29, 0.02059199907399667, 0.0003122853038535212),
    @(30, 0.01183140737501635, -0.01578616807981259),
    @(31, 0.008363753543070229, 0.015343083652269),
    @(32, 0.007740130312104565, -0.01048667011299009),
    @(33, 0.008529805374963471, 0.01486011057128489),
    @(34, 0.007152232518293211, -0.007232767232767223),
    @(35, $null, 0.003767461700544983)
)

foreach ($update in $updates) {
    $row = $update[0]
    $newWeight = $update[1]
    $newPctChange = $update[2]

    if ($null -ne $newWeight) {
        $ws.Cells.Item($row, 4).Value = $newWeight
    }
    $ws.Cells.Item($row, 5).Value = $newPctChange
}

if ($wasProtected) {
    $ws.Protect("D382")
}
